$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header note / timestamp update (A1) ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 22 de Mayo de 2020 a las 13:35"

# --- Country data refresh (newer snapshot) ---
# columns: A=Pais B=Casos totales C=Nuevos casos D=Casos activos
#          E=Recuperados F=Casos criticos G=Muertes hoy H=Muertes
# Table is ranked by "Casos totales" descending; several countries swap
# places with their neighbours as their numbers update.

function Set-Row($r, $pais, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($r, 1).Value = $pais
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
    $ws.Cells.Item($r, 6).Value = $f
    $ws.Cells.Item($r, 7).Value = $g
    $ws.Cells.Item($r, 8).Value = $h
}

# Alemania (row 11) - stats refresh only, no reordering
Set-Row 11 "Alemania" 179156 135 159000 11840 0 7 8316

# India (row 14) - stats refresh only, no reordering
Set-Row 14 "India" 119524 1298 49026 66898 0 16 3600

# Kuwait moves ahead of Sudafrica (rows 37/38 swap ranking)
Set-Row 37 "Kuwait" 19564 955 5515 13911 0 9 138
Set-Row 38 "Sudafrica" 19137 0 8950 9818 0 0 369

# Senegal moves ahead of Grecia (rows 78/79 swap ranking)
Set-Row 78 "Senegal" 2909 97 1311 1565 0 0 33
Set-Row 79 "Grecia" 2853 0 1374 1311 0 0 168

# Madagascar moves ahead of Taiwan, Republica de Africa Central, Etiopia,
# Estado de Palestina (rows 135-139 shift down one rank)
Set-Row 135 "Madagascar" 448 43 135 311 0 0 2
Set-Row 136 "Taiwan" 441 0 408 26 0 0 7
Set-Row 137 "Republica de Africa Central" 436 0 18 418 0 0 0
Set-Row 138 "Etiopia" 429 30 128 296 0 0 5
Set-Row 139 "Estado de Palestina" 423 0 346 75 0 0 2

# Groenlandia / Seychelles move ahead of Montserrat (rows 209-211 reorder)
Set-Row 209 "Groenlandia" 11 0 11 0 0 0 0
Set-Row 210 "Seychelles" 11 0 11 0 0 0 0
Set-Row 211 "Montserrat" 11 0 10 0 0 0 1

# Bonaire, San Eustaquio y Saba moves ahead of Sahara Occidental (rows 214/215 swap)
Set-Row 214 "Bonaire, San Eustaquio y Saba" 6 0 6 0 0 0 0
Set-Row 215 "Sahara Occidental" 6 0 6 0 0 0 0
